$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.238.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.428.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "406.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("E7").Value = "  -2.48%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.692"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.11%  "

$ws.Range("E10").Value = "  +7.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.421.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.236.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.59%  "

$ws.Range("E18").Value = "  -1.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000149"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.59%  "

$ws.Range("E20").Value = "  -3.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "83.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "311.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.89%  "

$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.53%  "

$ws.Range("E27").Value = "  -5.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "44.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.07%  "

$ws.Range("E31").Value = "  -1.21%  "

$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.06%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0485"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.323"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.73%  "

$ws.Range("E40").Value = "  -4.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "142.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.63%  "

$ws.Range("E42").Value = "  -0.48%  "

$ws.Range("E43").Value = "  -3.16%  "

$ws.Range("E44").Value = "  -3.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.64%  "

$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.101.28"
$ws.Range("D48").Style = "Normal"

$ws.Range("E49").Value = "  +2.85%  "

$ws.Range("E50").Value = "  -1.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +26.99%  "
